$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.512427745349363
$ws.Range("C2").Value = -1.363767227788882
$ws.Range("D2").Value = 59.70675209903689
$ws.Range("E2").Value = -0.0003999063564507079
$ws.Range("F2").Value = 0.6898044308789899
$ws.Range("G2").Value = -0.5004414916667799
$ws.Range("H2").Value = -0.6350004731026235
$ws.Range("I2").Value = 1.621664813878342
$ws.Range("J2").Value = 4.355285082868471
$ws.Range("K2").Value = 33
$ws.Range("L2").Value = -8.080333650738645
$ws.Range("N2").Value = 4.355285082866482
$ws.Range("O2").Value = 5.021951749533149

$ws.Range("B3").Value = 6.780578827585828
$ws.Range("C3").Value = 53.10279935400346
$ws.Range("D3").Value = -89.19799699470447
$ws.Range("E3").Value = -0.00005418064916170746
$ws.Range("F3").Value = 0.37163061194269
$ws.Range("G3").Value = -0.3368040805509991
$ws.Range("H3").Value = 0.0386712405137799
$ws.Range("I3").Value = 1.890841099281444
$ws.Range("J3").Value = 4.355295902663592
$ws.Range("K3").Value = 76
$ws.Range("L3").Value = 64.70060460020855
$ws.Range("N3").Value = 4.355295902660338
$ws.Range("O3").Value = 5.021962569327004

$ws.Range("B4").Value = -0.01513397410389549
$ws.Range("C4").Value = 25.15261976488516
$ws.Range("D4").Value = 44.35061839695146
$ws.Range("E4").Value = 0.2877657036915104
$ws.Range("F4").Value = 1.18443800501441
$ws.Range("G4").Value = -0.5605967116354871
$ws.Range("H4").Value = -0.9456454821981288
$ws.Range("I4").Value = 0.8338599972250136
$ws.Range("J4").Value = 4.355369250924604
$ws.Range("K4").Value = 36
$ws.Range("L4").Value = -5.761773861116499
$ws.Range("N4").Value = 4.355369250931052
$ws.Range("O4").Value = 5.022035917597719

$ws.Range("B5").Value = 5.379649467481933
$ws.Range("C5").Value = -2.599789483648555
$ws.Range("D5").Value = 20.55006959486832
$ws.Range("E5").Value = 50.21933496961683
$ws.Range("F5").Value = 0.5961330041681703
$ws.Range("G5").Value = 0.6765836970136485
$ws.Range("H5").Value = -0.3212336317668543
$ws.Range("I5").Value = -0.3538195067581045
$ws.Range("J5").Value = 4.354102017697073
$ws.Range("K5").Value = 65
$ws.Range("L5").Value = -33.0665014051261
$ws.Range("N5").Value = 4.355385754896673
$ws.Range("O5").Value = 5.02205242156334

$ws.Range("B6").Value = -18.14343379251174
$ws.Range("C6").Value = 49.29212575291513
$ws.Range("D6").Value = 1.855967475239631
$ws.Range("E6").Value = 31.11628356901486
$ws.Range("F6").Value = -0.6542334917483561
$ws.Range("G6").Value = -0.6422895521552101
$ws.Range("H6").Value = 0.489667206518646
$ws.Range("I6").Value = -0.3677553154227562
$ws.Range("J6").Value = 4.354322482717333
$ws.Range("K6").Value = 53
$ws.Range("L6").Value = -18.49407682031511
$ws.Range("N6").Value = 4.355388030443549
$ws.Range("O6").Value = 5.022054697110216

$ws.Range("B7").Value = -144.0501276637553
$ws.Range("C7").Value = 144.1985215536318
$ws.Range("D7").Value = 0.02834899075494345
$ws.Range("E7").Value = -184.1771139912038
$ws.Range("F7").Value = 0.006392393086723125
$ws.Range("G7").Value = -0.3762215925428569
$ws.Range("H7").Value = 1.032905282683497
$ws.Range("I7").Value = -0.1155008040902508
$ws.Range("J7").Value = 4.355421483559979
$ws.Range("K7").Value = 40
$ws.Range("L7").Value = 234.3095798628854
$ws.Range("N7").Value = 4.355421483560008
$ws.Range("O7").Value = 5.022088150226675

$ws.Range("B8").Value = 2.242833893024097
$ws.Range("C8").Value = 4.590921744530472
$ws.Range("D8").Value = 69.50350709815329
$ws.Range("E8").Value = -0.0002337489343013588
$ws.Range("F8").Value = -0.005662553358535494
$ws.Range("G8").Value = 0.3728013791281621
$ws.Range("H8").Value = -0.4526863167629085
$ws.Range("I8").Value = 1.483440877429489
$ws.Range("J8").Value = 4.353957525117924
$ws.Range("K8").Value = 45
$ws.Range("L8").Value = -29.58683169051387
$ws.Range("N8").Value = 4.355424018045394
$ws.Range("O8").Value = 5.022090684712061

$ws.Range("B9").Value = 182.6048700990076
$ws.Range("C9").Value = -214.5918842984699
$ws.Range("D9").Value = 0.05227083228840537
$ws.Range("E9").Value = 48.06839951286155
$ws.Range("F9").Value = -0.382876534643398
$ws.Range("G9").Value = -0.2128477500075963
$ws.Range("H9").Value = 0.9543968675431658
$ws.Range("I9").Value = -0.2398966450387012
$ws.Range("J9").Value = 4.355433165429133
$ws.Range("K9").Value = 75
$ws.Range("L9").Value = 35.59730982517511
$ws.Range("N9").Value = 4.355433165429506
$ws.Range("O9").Value = 5.022099832096173

$ws.Range("B10").Value = -48.60009858307789
$ws.Range("C10").Value = 100.7528450272593
$ws.Range("D10").Value = 3.078983553176508
$ws.Range("E10").Value = -0.0001434597757417982
$ws.Range("F10").Value = 0.09617584565183712
$ws.Range("G10").Value = -1.297390320175113
$ws.Range("H10").Value = 0.4772201815601758
$ws.Range("I10").Value = 1.522118501330372
$ws.Range("J10").Value = 4.355579008451969
$ws.Range("K10").Value = 29
$ws.Range("L10").Value = 54.41205601879853
$ws.Range("N10").Value = 4.355579008452018
$ws.Range("O10").Value = 5.022245675118685

$ws.Range("B11").Value = 8.153141052053776
$ws.Range("C11").Value = 75.97522598138652
$ws.Range("D11").Value = 0.001186883175047206
$ws.Range("E11").Value = 0.4235767403985576
$ws.Range("F11").Value = -0.6925446724403637
$ws.Range("G11").Value = -0.9513589618600715
$ws.Range("H11").Value = 1.361170704013164
$ws.Range("I11").Value = 0.6545411439990407
$ws.Range("J11").Value = 4.35567484706745
$ws.Range("K11").Value = 49
$ws.Range("L11").Value = -4.001917544265474
$ws.Range("N11").Value = 4.355674847067466
$ws.Range("O11").Value = 5.022341513734133
